$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 808.7379963333333
$ws.Cells.Item(2, 8).Value = 2426.213989
$ws.Cells.Item(2, 9).Value = 0.2987634987814484
$ws.Cells.Item(2, 10).Value = 0.2987634987814484
$ws.Cells.Item(2, 13).Value = 16.27546433333333
$ws.Cells.Item(2, 14).Value = 48.826393
$ws.Cells.Item(2, 15).Value = 0.06628560529319844
$ws.Cells.Item(2, 16).Value = 0.06628560529319844
$ws.Cells.Item(2, 17).Value = 13162.58641433463
$ws.Cells.Item(2, 18).Value = 118463.2777290117
$ws.Cells.Item(2, 19).Value = 0.01980371935624206
$ws.Cells.Item(2, 20).Value = 0.01980371935624206
$ws.Cells.Item(3, 7).Value = 808.7379963333333
$ws.Cells.Item(3, 8).Value = 2426.213989
$ws.Cells.Item(3, 9).Value = 0.2987634987814484
$ws.Cells.Item(3, 10).Value = 0.2987634987814484
$ws.Cells.Item(3, 15).Value = 0.3480686258826592
$ws.Cells.Item(3, 16).Value = 0.3480686258826592
$ws.Cells.Item(3, 17).Value = 69117.31960557835
$ws.Cells.Item(3, 18).Value = 622055.8764502051
$ws.Cells.Item(3, 19).Value = 0.1039902004847543
$ws.Cells.Item(3, 20).Value = 0.1039902004847543
$ws.Cells.Item(4, 7).Value = 808.7379963333333
$ws.Cells.Item(4, 8).Value = 2426.213989
$ws.Cells.Item(4, 9).Value = 0.2987634987814484
$ws.Cells.Item(4, 10).Value = 0.2987634987814484
$ws.Cells.Item(4, 13).Value = 42.61351133333333
$ws.Cells.Item(4, 14).Value = 127.840534
$ws.Cells.Item(4, 15).Value = 0.17355341356458
$ws.Cells.Item(4, 16).Value = 0.17355341356458
$ws.Cells.Item(4, 17).Value = 34463.16577244778
$ws.Cells.Item(4, 18).Value = 310168.4919520301
$ws.Cells.Item(4, 19).Value = 0.05185142506201761
$ws.Cells.Item(4, 20).Value = 0.05185142506201761
$ws.Cells.Item(5, 7).Value = 808.7379963333333
$ws.Cells.Item(5, 8).Value = 2426.213989
$ws.Cells.Item(5, 9).Value = 0.2987634987814484
$ws.Cells.Item(5, 10).Value = 0.2987634987814484
$ws.Cells.Item(5, 13).Value = 101.183272
$ws.Cells.Item(5, 14).Value = 303.549816
$ws.Cells.Item(5, 15).Value = 0.4120923552595624
$ws.Cells.Item(5, 16).Value = 0.4120923552595624
$ws.Cells.Item(5, 17).Value = 81830.75665973066
$ws.Cells.Item(5, 18).Value = 736476.8099375759
$ws.Cells.Item(5, 19).Value = 0.1231181538784345
$ws.Cells.Item(5, 20).Value = 0.1231181538784345
$ws.Cells.Item(6, 9).Value = 0.3611891679133081
$ws.Cells.Item(6, 10).Value = 0.361189167913308
$ws.Cells.Item(6, 13).Value = 16.27546433333333
$ws.Cells.Item(6, 14).Value = 48.826393
$ws.Cells.Item(6, 15).Value = 0.06628560529319844
$ws.Cells.Item(6, 16).Value = 0.06628560529319844
$ws.Cells.Item(6, 17).Value = 15912.86637748984
$ws.Cells.Item(6, 18).Value = 143215.7973974086
$ws.Cells.Item(6, 19).Value = 0.02394164262048031
$ws.Cells.Item(6, 20).Value = 0.02394164262048031
$ws.Cells.Item(7, 9).Value = 0.3611891679133081
$ws.Cells.Item(7, 10).Value = 0.361189167913308
$ws.Cells.Item(7, 15).Value = 0.3480686258826592
$ws.Cells.Item(7, 16).Value = 0.3480686258826592
$ws.Cells.Item(7, 19).Value = 0.1257186173592862
$ws.Cells.Item(7, 20).Value = 0.1257186173592862
$ws.Cells.Item(8, 9).Value = 0.3611891679133081
$ws.Cells.Item(8, 10).Value = 0.361189167913308
$ws.Cells.Item(8, 13).Value = 42.61351133333333
$ws.Cells.Item(8, 14).Value = 127.840534
$ws.Cells.Item(8, 15).Value = 0.17355341356458
$ws.Cells.Item(8, 16).Value = 0.17355341356458
$ws.Cells.Item(8, 17).Value = 41664.13306772317
$ws.Cells.Item(8, 18).Value = 374977.1976095085
$ws.Cells.Item(8, 19).Value = 0.06268561303390488
$ws.Cells.Item(8, 20).Value = 0.06268561303390488
$ws.Cells.Item(9, 9).Value = 0.3611891679133081
$ws.Cells.Item(9, 10).Value = 0.361189167913308
$ws.Cells.Item(9, 13).Value = 101.183272
$ws.Cells.Item(9, 14).Value = 303.549816
$ws.Cells.Item(9, 15).Value = 0.4120923552595624
$ws.Cells.Item(9, 16).Value = 0.4120923552595624
$ws.Cells.Item(9, 17).Value = 98929.02924284471
$ws.Cells.Item(9, 18).Value = 890361.2631856024
$ws.Cells.Item(9, 19).Value = 0.1488432948996367
$ws.Cells.Item(9, 20).Value = 0.1488432948996367
$ws.Cells.Item(10, 7).Value = 480.4688006666667
$ws.Cells.Item(10, 8).Value = 1441.406402
$ws.Cells.Item(10, 9).Value = 0.1774944921511204
$ws.Cells.Item(10, 10).Value = 0.1774944921511204
$ws.Cells.Item(10, 13).Value = 16.27546433333333
$ws.Cells.Item(10, 14).Value = 48.826393
$ws.Cells.Item(10, 15).Value = 0.06628560529319844
$ws.Cells.Item(10, 16).Value = 0.06628560529319844
$ws.Cells.Item(10, 17).Value = 7819.852828529776
$ws.Cells.Item(10, 18).Value = 70378.67545676799
$ws.Cells.Item(10, 19).Value = 0.01176532984844587
$ws.Cells.Item(10, 20).Value = 0.01176532984844587
$ws.Cells.Item(11, 7).Value = 480.4688006666667
$ws.Cells.Item(11, 8).Value = 1441.406402
$ws.Cells.Item(11, 9).Value = 0.1774944921511204
$ws.Cells.Item(11, 10).Value = 0.1774944921511204
$ws.Cells.Item(11, 15).Value = 0.3480686258826592
$ws.Cells.Item(11, 16).Value = 0.3480686258826592
$ws.Cells.Item(11, 17).Value = 41062.39079497813
$ws.Cells.Item(11, 18).Value = 369561.5171548032
$ws.Cells.Item(11, 19).Value = 0.06178026398478091
$ws.Cells.Item(11, 20).Value = 0.0617802639847809
$ws.Cells.Item(12, 7).Value = 480.4688006666667
$ws.Cells.Item(12, 8).Value = 1441.406402
$ws.Cells.Item(12, 9).Value = 0.1774944921511204
$ws.Cells.Item(12, 10).Value = 0.1774944921511204
$ws.Cells.Item(12, 13).Value = 42.61351133333333
$ws.Cells.Item(12, 14).Value = 127.840534
$ws.Cells.Item(12, 15).Value = 0.17355341356458
$ws.Cells.Item(12, 16).Value = 0.17355341356458
$ws.Cells.Item(12, 17).Value = 20474.46268252207
$ws.Cells.Item(12, 18).Value = 184270.1641426987
$ws.Cells.Item(12, 19).Value = 0.03080477500173849
$ws.Cells.Item(12, 20).Value = 0.03080477500173849
$ws.Cells.Item(13, 7).Value = 480.4688006666667
$ws.Cells.Item(13, 8).Value = 1441.406402
$ws.Cells.Item(13, 9).Value = 0.1774944921511204
$ws.Cells.Item(13, 10).Value = 0.1774944921511204
$ws.Cells.Item(13, 13).Value = 101.183272
$ws.Cells.Item(13, 14).Value = 303.549816
$ws.Cells.Item(13, 15).Value = 0.4120923552595624
$ws.Cells.Item(13, 16).Value = 0.4120923552595624
$ws.Cells.Item(13, 17).Value = 48615.40534536912
$ws.Cells.Item(13, 18).Value = 437538.648108322
$ws.Cells.Item(13, 19).Value = 0.07314412331615512
$ws.Cells.Item(13, 20).Value = 0.07314412331615511
$ws.Cells.Item(14, 7).Value = 440.0224913333333
$ws.Cells.Item(14, 8).Value = 1320.067474
$ws.Cells.Item(14, 9).Value = 0.1625528411541232
$ws.Cells.Item(14, 10).Value = 0.1625528411541232
$ws.Cells.Item(14, 13).Value = 16.27546433333333
$ws.Cells.Item(14, 14).Value = 48.826393
$ws.Cells.Item(14, 15).Value = 0.06628560529319844
$ws.Cells.Item(14, 16).Value = 0.06628560529319844
$ws.Cells.Item(14, 17).Value = 7161.570363560142
$ws.Cells.Item(14, 18).Value = 64454.13327204128
$ws.Cells.Item(14, 19).Value = 0.0107749134680302
$ws.Cells.Item(14, 20).Value = 0.0107749134680302
$ws.Cells.Item(15, 7).Value = 440.0224913333333
$ws.Cells.Item(15, 8).Value = 1320.067474
$ws.Cells.Item(15, 9).Value = 0.1625528411541232
$ws.Cells.Item(15, 10).Value = 0.1625528411541232
$ws.Cells.Item(15, 15).Value = 0.3480686258826592
$ws.Cells.Item(15, 16).Value = 0.3480686258826592
$ws.Cells.Item(15, 17).Value = 37605.7206474983
$ws.Cells.Item(15, 18).Value = 338451.4858274846
$ws.Cells.Item(15, 19).Value = 0.05657954405383785
$ws.Cells.Item(15, 20).Value = 0.05657954405383784
$ws.Cells.Item(16, 7).Value = 440.0224913333333
$ws.Cells.Item(16, 8).Value = 1320.067474
$ws.Cells.Item(16, 9).Value = 0.1625528411541232
$ws.Cells.Item(16, 10).Value = 0.1625528411541232
$ws.Cells.Item(16, 13).Value = 42.61351133333333
$ws.Cells.Item(16, 14).Value = 127.840534
$ws.Cells.Item(16, 15).Value = 0.17355341356458
$ws.Cells.Item(16, 16).Value = 0.17355341356458
$ws.Cells.Item(16, 17).Value = 18750.90342135457
$ws.Cells.Item(16, 18).Value = 168758.1307921911
$ws.Cells.Item(16, 19).Value = 0.02821160046691903
$ws.Cells.Item(16, 20).Value = 0.02821160046691903
$ws.Cells.Item(17, 7).Value = 440.0224913333333
$ws.Cells.Item(17, 8).Value = 1320.067474
$ws.Cells.Item(17, 9).Value = 0.1625528411541232
$ws.Cells.Item(17, 10).Value = 0.1625528411541232
$ws.Cells.Item(17, 13).Value = 101.183272
$ws.Cells.Item(17, 14).Value = 303.549816
$ws.Cells.Item(17, 15).Value = 0.4120923552595624
$ws.Cells.Item(17, 16).Value = 0.4120923552595624
$ws.Cells.Item(17, 17).Value = 44522.9154266983
$ws.Cells.Item(17, 18).Value = 400706.2388402847
$ws.Cells.Item(17, 19).Value = 0.06698678316533617
$ws.Cells.Item(17, 20).Value = 0.06698678316533616
